$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 70719.71000000001
$ws.Range("J21").Value = 59004.2
$ws.Range("L21").Value = 59004.2
$ws.Range("N21").Value = -59940.2
$ws.Range("H23").Value = 70719.71000000001
$ws.Range("J23").Value = 59004.2
$ws.Range("L23").Value = 59004.2
$ws.Range("N23").Value = -59472.2
$ws.Range("H64").Value = 6999.7393
$ws.Range("J64").Value = 7399.9
$ws.Range("L64").Value = 7399.9
$ws.Range("N64").Value = -7895.9
$ws.Range("H67").Value = 6999.7393
$ws.Range("J67").Value = 7399.9
$ws.Range("L67").Value = 7399.9
$ws.Range("N67").Value = -9115.9
$ws.Range("H69").Value = 10014.964
$ws.Range("H72").Value = 10014.964
$ws.Range("H112").Value = 1044154.25
$ws.Range("J112").Value = 1451735.6
$ws.Range("L112").Value = 4355206.800000001
$ws.Range("N112").Value = -4357422.800000001
$ws.Range("H113").Value = 99999.5
$ws.Range("I113").Value = 99999.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 99999.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -96745.5
$ws.Range("N113").ClearContents()
$ws.Range("H125").Value = 3778.875
$ws.Range("I125").Value = 3949
$ws.Range("J125").Value = 3754.5715
$ws.Range("K125").Value = 35541
$ws.Range("L125").Value = 33791.1435
$ws.Range("M125").Value = -33081
$ws.Range("N125").Value = -38711.1435
$ws.Range("H129").Value = 1368
$ws.Range("I129").Value = 654.8570999999999
$ws.Range("K129").Value = 1964.5713
$ws.Range("M129").Value = 3035.4287
$ws.Range("H138").Value = 2562.44
$ws.Range("I138").Value = 1213.909
$ws.Range("K138").Value = 3641.727
$ws.Range("M138").Value = 1498.273

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3067.625
$ws.Range("I2").Value = 3091.625
$ws.Range("J2").Value = 3043.625
$ws.Range("K2").Value = 3091.625
$ws.Range("L2").Value = 3043.625
$ws.Range("M2").Value = -2978.625
$ws.Range("N2").Value = -3269.625
$ws.Range("H45").Value = 4895.6924
$ws.Range("I45").Value = 4331.273
$ws.Range("K45").Value = 4331.273
$ws.Range("M45").Value = -3954.273
$ws.Range("H63").Value = 3856.9
$ws.Range("I63").Value = 2400
$ws.Range("K63").Value = 2400
$ws.Range("M63").Value = -1714
$ws.Range("H66").Value = 3856.9
$ws.Range("I66").Value = 2400
$ws.Range("K66").Value = 12000
$ws.Range("M66").Value = -8568
$ws.Range("H116").Value = 3067.625
$ws.Range("I116").Value = 3091.625
$ws.Range("J116").Value = 3043.625
$ws.Range("K116").Value = 3091.625
$ws.Range("L116").Value = 3043.625
$ws.Range("M116").Value = -797.625
$ws.Range("N116").Value = -7631.625
$ws.Range("H122").Value = 3145
$ws.Range("I122").Value = 3074.375
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 9223.125
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -6773.125
$ws.Range("N122").Value = -14899.9999
$ws.Range("H131").Value = 125657.5
$ws.Range("J131").Value = 125657.5
$ws.Range("L131").Value = 125657.5
$ws.Range("N131").Value = -135737.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3067.625
$ws.Range("I3").Value = 3091.625
$ws.Range("J3").Value = 3043.625
$ws.Range("K3").Value = 3091.625
$ws.Range("L3").Value = 3043.625
$ws.Range("M3").Value = -2977.625
$ws.Range("N3").Value = -3271.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 455702.28
$ws.Range("I132").Value = 599.8182
$ws.Range("J132").Value = 910804.75
$ws.Range("K132").Value = 5398.3638
$ws.Range("L132").Value = 8197242.75
$ws.Range("M132").Value = -2868.3638
$ws.Range("N132").Value = -8202302.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 50001
$ws.Range("J123").Value = 50001
$ws.Range("L123").Value = 50001
$ws.Range("N123").Value = -54901
$ws.Range("H132").Value = 3389.647
$ws.Range("I132").Value = 3098.6924
$ws.Range("J132").Value = 4335.25
$ws.Range("K132").Value = 9296.0772
$ws.Range("L132").Value = 13005.75
$ws.Range("M132").Value = -6766.0772
$ws.Range("N132").Value = -18065.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10298.667
$ws.Range("I7").Value = 10666.667
$ws.Range("J7").Value = 9930.666999999999
$ws.Range("K7").Value = 10666.667
$ws.Range("L7").Value = 9930.666999999999
$ws.Range("M7").Value = -10554.667
$ws.Range("N7").Value = -10154.667
$ws.Range("H26").Value = 20009
$ws.Range("I26").Value = 20009
$ws.Range("K26").Value = 20009
$ws.Range("M26").Value = -19714
$ws.Range("H93").Value = 3161
$ws.Range("I93").Value = 2951.25
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 2951.25
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -1703.25
$ws.Range("N93").Value = -6496
$ws.Range("H122").Value = 30275.572
$ws.Range("I122").Value = 34235.5
$ws.Range("K122").Value = 102706.5
$ws.Range("M122").Value = -100256.5
$ws.Range("H126").Value = 10298.667
$ws.Range("I126").Value = 10666.667
$ws.Range("J126").Value = 9930.666999999999
$ws.Range("K126").Value = 32000.001
$ws.Range("L126").Value = 29792.001
$ws.Range("M126").Value = -29530.001
$ws.Range("N126").Value = -34732.001
$ws.Range("H132").Value = 4246.485
$ws.Range("I132").Value = 3829.36
$ws.Range("K132").Value = 11488.08
$ws.Range("M132").Value = -8958.08

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2464.8845
$ws.Range("I81").Value = 1837.5714
$ws.Range("K81").Value = 3675.1428
$ws.Range("M81").Value = -2614.1428
$ws.Range("H84").Value = 2464.8845
$ws.Range("I84").Value = 1837.5714
$ws.Range("K84").Value = 18375.714
$ws.Range("M84").Value = -13071.714
$ws.Range("H130").Value = 56872.4
$ws.Range("J130").Value = 56872.4
$ws.Range("L130").Value = 56872.4
$ws.Range("N130").Value = -66912.39999999999
$ws.Range("H132").Value = 1925.8334
$ws.Range("I132").Value = 1555.4546
$ws.Range("K132").Value = 4666.3638
$ws.Range("M132").Value = -2136.3638
